$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.653906
$ws.Range("H2").Value = 1.961718
$ws.Range("I2").Value = 0.00670030715761011
$ws.Range("J2").Value = 0.00670030715761011
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 30.99161333333333
$ws.Range("N2").Value = 92.97484
$ws.Range("O2").Value = 0.3599121977633812
$ws.Range("P2").Value = 0.3599121977633811
$ws.Range("Q2").Value = 20.26560190834667
$ws.Range("R2").Value = 182.39041717512
$ws.Range("S2").Value = 0.002411522274785169
$ws.Range("T2").Value = 0.002411522274785168

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.653906
$ws.Range("H3").Value = 1.961718
$ws.Range("I3").Value = 0.00670030715761011
$ws.Range("J3").Value = 0.00670030715761011
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 29.913269
$ws.Range("N3").Value = 89.739807
$ws.Range("O3").Value = 0.3473891556493311
$ws.Range("P3").Value = 0.3473891556493311
$ws.Range("Q3").Value = 19.560466078714
$ws.Range("R3").Value = 176.044194708426
$ws.Range("S3").Value = 0.002327614046073346
$ws.Range("T3").Value = 0.002327614046073346

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.653906
$ws.Range("H4").Value = 1.961718
$ws.Range("I4").Value = 0.00670030715761011
$ws.Range("J4").Value = 0.00670030715761011
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 19.150218
$ws.Range("N4").Value = 57.450654
$ws.Range("O4").Value = 0.2223955550134164
$ws.Range("P4").Value = 0.2223955550134163
$ws.Range("Q4").Value = 12.522442451508
$ws.Range("R4").Value = 112.701982063572
$ws.Range("S4").Value = 0.001490118529077067
$ws.Range("T4").Value = 0.001490118529077067

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.653906
$ws.Range("H5").Value = 1.961718
$ws.Range("I5").Value = 0.00670030715761011
$ws.Range("J5").Value = 0.00670030715761011
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 6.053716000000001
$ws.Range("N5").Value = 18.161148
$ws.Range("O5").Value = 0.07030309157387134
$ws.Range("P5").Value = 0.07030309157387132
$ws.Range("Q5").Value = 3.958561214696
$ws.Range("R5").Value = 35.627050932264
$ws.Range("S5").Value = 0.0004710523076745291
$ws.Range("T5").Value = 0.000471052307674529

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 96.11977900000001
$ws.Range("H6").Value = 288.359337
$ws.Range("I6").Value = 0.984900036429704
$ws.Range("J6").Value = 0.984900036429704
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 30.99161333333333
$ws.Range("N6").Value = 92.97484
$ws.Range("O6").Value = 0.3599121977633812
$ws.Range("P6").Value = 0.3599121977633811
$ws.Range("Q6").Value = 2978.907024453454
$ws.Range("R6").Value = 26810.16322008109
$ws.Range("S6").Value = 0.354477536688649
$ws.Range("T6").Value = 0.3544775366886489

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 96.11977900000001
$ws.Range("H7").Value = 288.359337
$ws.Range("I7").Value = 0.984900036429704
$ws.Range("J7").Value = 0.984900036429704
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 29.913269
$ws.Range("N7").Value = 89.739807
$ws.Range("O7").Value = 0.3473891556493311
$ws.Range("P7").Value = 0.3473891556493311
$ws.Range("Q7").Value = 2875.256805447551
$ws.Range("R7").Value = 25877.31124902796
$ws.Range("S7").Value = 0.3421435920543104
$ws.Range("T7").Value = 0.3421435920543103

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 96.11977900000001
$ws.Range("H8").Value = 288.359337
$ws.Range("I8").Value = 0.984900036429704
$ws.Range("J8").Value = 0.984900036429704
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 19.150218
$ws.Range("N8").Value = 57.450654
$ws.Range("O8").Value = 0.2223955550134164
$ws.Range("P8").Value = 0.2223955550134163
$ws.Range("Q8").Value = 1840.714721961822
$ws.Range("R8").Value = 16566.4324976564
$ws.Range("S8").Value = 0.2190373902345181
$ws.Range("T8").Value = 0.219037390234518

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 96.11977900000001
$ws.Range("H9").Value = 288.359337
$ws.Range("I9").Value = 0.984900036429704
$ws.Range("J9").Value = 0.984900036429704
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 6.053716000000001
$ws.Range("N9").Value = 18.161148
$ws.Range("O9").Value = 0.07030309157387134
$ws.Range("P9").Value = 0.07030309157387132
$ws.Range("Q9").Value = 581.8818440487642
$ws.Range("R9").Value = 5236.936596438877
$ws.Range("S9").Value = 0.0692415174522267
$ws.Range("T9").Value = 0.06924151745222669

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.2341223333333333
$ws.Range("H10").Value = 0.702367
$ws.Range("I10").Value = 0.002398955730318598
$ws.Range("J10").Value = 0.002398955730318598
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 30.99161333333333
$ws.Range("N10").Value = 92.97484
$ws.Range("O10").Value = 0.3599121977633812
$ws.Range("P10").Value = 0.3599121977633811
$ws.Range("Q10").Value = 7.255828827364444
$ws.Range("R10").Value = 65.30245944628
$ws.Range("S10").Value = 0.0008634134292360241
$ws.Range("T10").Value = 0.0008634134292360236

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.2341223333333333
$ws.Range("H11").Value = 0.702367
$ws.Range("I11").Value = 0.002398955730318598
$ws.Range("J11").Value = 0.002398955730318598
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 29.913269
$ws.Range("N11").Value = 89.739807
$ws.Range("O11").Value = 0.3473891556493311
$ws.Range("P11").Value = 0.3473891556493311
$ws.Range("Q11").Value = 7.003364335907666
$ws.Range("R11").Value = 63.03027902316899
$ws.Range("S11").Value = 0.0008333712055955025
$ws.Range("T11").Value = 0.0008333712055955022

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.2341223333333333
$ws.Range("H12").Value = 0.702367
$ws.Range("I12").Value = 0.002398955730318598
$ws.Range("J12").Value = 0.002398955730318598
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 19.150218
$ws.Range("N12").Value = 57.450654
$ws.Range("O12").Value = 0.2223955550134164
$ws.Range("P12").Value = 0.2223955550134163
$ws.Range("Q12").Value = 4.483493722002
$ws.Range("R12").Value = 40.351443498018
$ws.Range("S12").Value = 0.0005335170910968204
$ws.Range("T12").Value = 0.0005335170910968201

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.2341223333333333
$ws.Range("H13").Value = 0.702367
$ws.Range("I13").Value = 0.002398955730318598
$ws.Range("J13").Value = 0.002398955730318598
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 6.053716000000001
$ws.Range("N13").Value = 18.161148
$ws.Range("O13").Value = 0.07030309157387134
$ws.Range("P13").Value = 0.07030309157387132
$ws.Range("Q13").Value = 1.417310115257333
$ws.Range("R13").Value = 12.755791037316
$ws.Range("S13").Value = 0.0001686540043902518
$ws.Range("T13").Value = 0.0001686540043902518

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.5856290000000001
$ws.Range("H14").Value = 1.756887
$ws.Range("I14").Value = 0.00600070068236727
$ws.Range("J14").Value = 0.006000700682367269
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 30.99161333333333
$ws.Range("N14").Value = 92.97484
$ws.Range("O14").Value = 0.3599121977633812
$ws.Range("P14").Value = 0.3599121977633811
$ws.Range("Q14").Value = 18.14958752478667
$ws.Range("R14").Value = 163.34628772308
$ws.Range("S14").Value = 0.002159725370711026
$ws.Range("T14").Value = 0.002159725370711025

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.5856290000000001
$ws.Range("H15").Value = 1.756887
$ws.Range("I15").Value = 0.00600070068236727
$ws.Range("J15").Value = 0.006000700682367269
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 29.913269
$ws.Range("N15").Value = 89.739807
$ws.Range("O15").Value = 0.3473891556493311
$ws.Range("P15").Value = 0.3473891556493311
$ws.Range("Q15").Value = 17.518077811201
$ws.Range("R15").Value = 157.662700300809
$ws.Range("S15").Value = 0.002084578343351931
$ws.Range("T15").Value = 0.00208457834335193

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.5856290000000001
$ws.Range("H16").Value = 1.756887
$ws.Range("I16").Value = 0.00600070068236727
$ws.Range("J16").Value = 0.006000700682367269
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 19.150218
$ws.Range("N16").Value = 57.450654
$ws.Range("O16").Value = 0.2223955550134164
$ws.Range("P16").Value = 0.2223955550134163
$ws.Range("Q16").Value = 11.214923017122
$ws.Range("R16").Value = 100.934307154098
$ws.Range("S16").Value = 0.001334529158724456
$ws.Range("T16").Value = 0.001334529158724455

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.5856290000000001
$ws.Range("H17").Value = 1.756887
$ws.Range("I17").Value = 0.00600070068236727
$ws.Range("J17").Value = 0.006000700682367269
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 6.053716000000001
$ws.Range("N17").Value = 18.161148
$ws.Range("O17").Value = 0.07030309157387134
$ws.Range("P17").Value = 0.07030309157387132
$ws.Range("Q17").Value = 3.545231647364001
$ws.Range("R17").Value = 31.907084826276
$ws.Range("S17").Value = 0.0004218678095798584
$ws.Range("T17").Value = 0.0004218678095798583
